# Minor tweaks to simplify age bias correction
# Update the SAD-by-age correction values in column B (rows 2-37) to the
# refreshed figures, then restore the active selection to B1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.096476606000000006
$ws.Range("B3").Value = 0.082428554000000001
$ws.Range("B4").Value = 0.071943749000000001
$ws.Range("B5").Value = 0.063603134000000006
$ws.Range("B6").Value = 0.056778598
$ws.Range("B7").Value = 0.051076533
$ws.Range("B8").Value = 0.046236073000000003
$ws.Range("B9").Value = 0.042074571999999998
$ws.Range("B10").Value = 0.038452368000000001
$ws.Range("B11").Value = 0.035252987
$ws.Range("B12").Value = 0.032379460999999998
$ws.Range("B13").Value = 0.029757838000000002
$ws.Range("B14").Value = 0.027346947
$ws.Range("B15").Value = 0.025158990999999999
$ws.Range("B16").Value = 0.023276809999999998
$ws.Range("B17").Value = 0.02169455
$ws.Range("B18").Value = 0.020293176
$ws.Range("B19").Value = 0.019023919
$ws.Range("B20").Value = 0.017946446000000001
$ws.Range("B21").Value = 0.023286067000000001
$ws.Range("B22").Value = 0.021274210000000002
$ws.Range("B23").Value = 0.019445661
$ws.Range("B24").Value = 0.017761559
$ws.Range("B25").Value = 0.016203492999999999
$ws.Range("B26").Value = 0.014755971
$ws.Range("B27").Value = 0.013405158
$ws.Range("B28").Value = 0.012137317
$ws.Range("B29").Value = 0.010936708
$ws.Range("B30").Value = 0.0097851799999999992
$ws.Range("B31").Value = 0.0086647670000000003
$ws.Range("B32").Value = 0.0075615320000000001
$ws.Range("B33").Value = 0.0064687579999999998
$ws.Range("B34").Value = 0.005392099
$ws.Range("B35").Value = 0.0043530590000000003
$ws.Range("B36").Value = 0.0033637049999999998
$ws.Range("B37").Value = 0.0040034440000000001

$ws.Range("B1").Select()

